# (#232) disabled failing cts tests - removed bad data from excel data providers
# - remove the two "bad data" rows from the pages_with_breadcrumbs sheet:
#     row 14: about-nci/organization/dcb/division-of-cancer-biology-test (Mini Landing Page)
#     row 18: news-events/press-releases/2019/brca-exchange-test (Press Release)
# - leave the other sheet's data untouched
# - the workbook now ends up with "pages_with_no_breadcrumbs" as the active/selected sheet,
#   and the last data row (now row 17) selected on "pages_with_breadcrumbs"

$wb = $excel.ActiveWorkbook

$wsBreadcrumbs = $wb.Worksheets.Item("pages_with_breadcrumbs")
$wsNoBreadcrumbs = $wb.Worksheets.Item("pages_with_no_breadcrumbs")

# Delete the higher-numbered row first so the row 14 index is not shifted
# before it is removed.
$wsBreadcrumbs.Rows.Item(18).Delete() | Out-Null
$wsBreadcrumbs.Rows.Item(14).Delete() | Out-Null

# Select the new last row (entire row) on the breadcrumbs sheet, matching
# the post-edit selection left behind by the deletion.
$wsBreadcrumbs.Range("A17:XFD17").Select() | Out-Null

# Make "pages_with_no_breadcrumbs" the active sheet/tab.
$wsNoBreadcrumbs.Activate() | Out-Null
